# Applies the "GIS foto Pekla zplozenci" edit:
#  - swaps the B1/C1 header labels (longitude/latitude -> latitude/longitude)
#  - adds a new "stav" (status) column F with header in F1 and "ok" in F2:F16
#  - adds a note "dnes Štěpánská" in E2
#  - moves the active selection to F16

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the longitude/latitude header labels in B1 and C1.
$ws.Range("B1").Value = "latitude"
$ws.Range("C1").Value = "longitude"

# New description note for row 2 (added to the shared-string table before
# the "stav"/"ok" strings, matching the saved workbook's string order).
$ws.Range("E2").Value = "dnes Štěpánská"

# New status column header.
$ws.Range("F1").Value = "stav"

# Fill "ok" status for rows 2 through 16 in the new column F.
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 6).Value = "ok"
}

# Update the visible selection / scroll position to match the saved file.
$ws.Range("F16").Select()
